$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Insert a new row at position 4; this shifts the existing row 4
# ("簡単なHP作成とAWS構築をしてくれる方募集") down to row 5, carrying its
# values and styles (including the Hyperlink cell style on column F) along
# with it.
$ws.Rows.Item(4).Insert()

# Refresh the "取得日時" timestamp on the two rows that were already present.
$ws.Range("A2").Value = "2026-01-04 18:25:59"
$ws.Range("A3").Value = "2026-01-04 18:25:59"

# Populate the newly inserted row 4 with the new listing.
$ws.Range("A4").Value = "2026-01-04 18:25:59"
$ws.Range("B4").Value = "【急募】kintone案件管理アプリにExcel見積計算式組込"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5464763"
$ws.Range("G4").Value = 55
$ws.Range("H4").Value = "◇アプリ"

# Refresh the timestamp on row 5 (the shifted former row 4).
$ws.Range("A5").Value = "2026-01-04 18:25:59"

# The row insert moves cell values/styles down but does not relocate the
# hyperlink relationships themselves, so rebuild the hyperlinks collection
# cleanly, in row order, for all four URL cells.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5464587")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5464763")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5457524")

# Make sure every URL cell keeps the shared "Hyperlink" cell style.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
